$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1.0
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.04979866666666666
$ws.Cells.Item(2, 8).Value = 0.149396
$ws.Cells.Item(2, 9).Value = 0.04800668898255549
$ws.Cells.Item(2, 10).Value = 0.04800668898255549
$ws.Cells.Item(2, 11).Value = 3.0
$ws.Cells.Item(2, 13).Value = 47.50824466666666
$ws.Cells.Item(2, 14).Value = 142.524734
$ws.Cells.Item(2, 15).Value = 0.04546113442798697
$ws.Cells.Item(2, 16).Value = 0.04619248897260801
$ws.Cells.Item(2, 17).Value = 2.365847240073777
$ws.Cells.Item(2, 18).Value = 21.292625160664
$ws.Cells.Item(2, 19).Value = 0.002182438541278516
$ws.Cells.Item(2, 20).Value = 0.002217548451438117
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1.0
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.04979866666666666
$ws.Cells.Item(3, 8).Value = 0.149396
$ws.Cells.Item(3, 9).Value = 0.04800668898255549
$ws.Cells.Item(3, 10).Value = 0.04800668898255549
$ws.Cells.Item(3, 11).Value = 3.0
$ws.Cells.Item(3, 15).Value = 0.105812544913079
$ws.Cells.Item(3, 16).Value = 0.1075148008416609
$ws.Cells.Item(3, 17).Value = 5.506600759035999
$ws.Cells.Item(3, 18).Value = 49.559406831324
$ws.Cells.Item(3, 19).Value = 0.005079709934094869
$ws.Cells.Item(3, 20).Value = 0.005161429605027011
$ws.Cells.Item(4, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4, 5).Value = 1.0
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.04979866666666666
$ws.Cells.Item(4, 8).Value = 0.149396
$ws.Cells.Item(4, 9).Value = 0.04800668898255549
$ws.Cells.Item(4, 10).Value = 0.04800668898255549
$ws.Cells.Item(4, 11).Value = 3.0
$ws.Cells.Item(4, 13).Value = 351.3736063333333
$ws.Cells.Item(4, 14).Value = 1054.120819
$ws.Cells.Item(4, 15).Value = 0.3362330657350935
$ws.Cells.Item(4, 16).Value = 0.3416422044152282
$ws.Cells.Item(4, 17).Value = 17.49793709725822
$ws.Cells.Item(4, 18).Value = 157.481433875324
$ws.Cells.Item(4, 19).Value = 0.01614143621239577
$ws.Cells.Item(4, 20).Value = 0.01640111105067651
$ws.Cells.Item(5, 4).Value = "MuSCs"
$ws.Cells.Item(5, 5).Value = 1.0
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.04979866666666666
$ws.Cells.Item(5, 8).Value = 0.149396
$ws.Cells.Item(5, 9).Value = 0.04800668898255549
$ws.Cells.Item(5, 10).Value = 0.04800668898255549
$ws.Cells.Item(5, 11).Value = 2.0
$ws.Cells.Item(5, 13).Value = 49.6371195
$ws.Cells.Item(5, 14).Value = 99.274239
$ws.Cells.Item(5, 15).Value = 0.04749827694204053
$ws.Cells.Item(5, 16).Value = 0.03217493596775667
$ws.Cells.Item(5, 17).Value = 2.471862368274
$ws.Cells.Item(5, 18).Value = 14.831174209644
$ws.Cells.Item(5, 19).Value = 0.002280235008363827
$ws.Cells.Item(5, 20).Value = 0.001544612144037733
$ws.Cells.Item(6, 5).Value = 1.0
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.04979866666666666
$ws.Cells.Item(6, 8).Value = 0.149396
$ws.Cells.Item(6, 9).Value = 0.04800668898255549
$ws.Cells.Item(6, 10).Value = 0.04800668898255549
$ws.Cells.Item(6, 11).Value = 3.0
$ws.Cells.Item(6, 13).Value = 485.9336543333334
$ws.Cells.Item(6, 14).Value = 1457.800963
$ws.Cells.Item(6, 15).Value = 0.4649949779817998
$ws.Cells.Item(6, 16).Value = 0.4724755698027463
$ws.Cells.Item(6, 17).Value = 24.19884807426089
$ws.Cells.Item(6, 18).Value = 217.789632668348
$ws.Cells.Item(6, 19).Value = 0.0223228692864225
$ws.Cells.Item(6, 20).Value = 0.02268198773137613
$ws.Cells.Item(7, 4).Value = "ECs"
$ws.Cells.Item(7, 5).Value = 3.0
$ws.Cells.Item(7, 7).Value = 0.809644
$ws.Cells.Item(7, 8).Value = 2.428932
$ws.Cells.Item(7, 9).Value = 0.7805094050963647
$ws.Cells.Item(7, 10).Value = 0.7805094050963647
$ws.Cells.Item(7, 11).Value = 3.0
$ws.Cells.Item(7, 13).Value = 47.50824466666666
$ws.Cells.Item(7, 14).Value = 142.524734
$ws.Cells.Item(7, 15).Value = 0.04546113442798697
$ws.Cells.Item(7, 16).Value = 0.04619248897260801
$ws.Cells.Item(7, 17).Value = 38.46476524489866
$ws.Cells.Item(7, 18).Value = 346.182887204088
$ws.Cells.Item(7, 19).Value = 0.03548284298739397
$ws.Cells.Item(7, 20).Value = 0.03605367208793066
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 3.0
$ws.Cells.Item(8, 7).Value = 0.809644
$ws.Cells.Item(8, 8).Value = 2.428932
$ws.Cells.Item(8, 9).Value = 0.7805094050963647
$ws.Cells.Item(8, 10).Value = 0.7805094050963647
$ws.Cells.Item(8, 11).Value = 3.0
$ws.Cells.Item(8, 15).Value = 0.105812544913079
$ws.Cells.Item(8, 16).Value = 0.1075148008416609
$ws.Cells.Item(8, 17).Value = 89.528225620812
$ws.Cells.Item(8, 18).Value = 805.754030587308
$ws.Cells.Item(8, 19).Value = 0.08258768648183967
$ws.Cells.Item(8, 20).Value = 0.08391631324397889
$ws.Cells.Item(9, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(9, 5).Value = 3.0
$ws.Cells.Item(9, 7).Value = 0.809644
$ws.Cells.Item(9, 8).Value = 2.428932
$ws.Cells.Item(9, 9).Value = 0.7805094050963647
$ws.Cells.Item(9, 10).Value = 0.7805094050963647
$ws.Cells.Item(9, 11).Value = 3.0
$ws.Cells.Item(9, 13).Value = 351.3736063333333
$ws.Cells.Item(9, 14).Value = 1054.120819
$ws.Cells.Item(9, 15).Value = 0.3362330657350935
$ws.Cells.Item(9, 16).Value = 0.3416422044152282
$ws.Cells.Item(9, 17).Value = 284.4875321261454
$ws.Cells.Item(9, 18).Value = 2560.387789135308
$ws.Cells.Item(9, 19).Value = 0.2624330701106247
$ws.Cells.Item(9, 20).Value = 0.2666549537239404
$ws.Cells.Item(10, 4).Value = "MuSCs"
$ws.Cells.Item(10, 5).Value = 3.0
$ws.Cells.Item(10, 7).Value = 0.809644
$ws.Cells.Item(10, 8).Value = 2.428932
$ws.Cells.Item(10, 9).Value = 0.7805094050963647
$ws.Cells.Item(10, 10).Value = 0.7805094050963647
$ws.Cells.Item(10, 11).Value = 2.0
$ws.Cells.Item(10, 13).Value = 49.6371195
$ws.Cells.Item(10, 14).Value = 99.274239
$ws.Cells.Item(10, 15).Value = 0.04749827694204053
$ws.Cells.Item(10, 16).Value = 0.03217493596775667
$ws.Cells.Item(10, 17).Value = 40.188395980458
$ws.Cells.Item(10, 18).Value = 241.130375882748
$ws.Cells.Item(10, 19).Value = 0.03707285187913443
$ws.Cells.Item(10, 20).Value = 0.02511284013120738
$ws.Cells.Item(11, 7).Value = 0.809644
$ws.Cells.Item(11, 8).Value = 2.428932
$ws.Cells.Item(11, 9).Value = 0.7805094050963647
$ws.Cells.Item(11, 10).Value = 0.7805094050963647
$ws.Cells.Item(11, 11).Value = 3.0
$ws.Cells.Item(11, 13).Value = 485.9336543333334
$ws.Cells.Item(11, 14).Value = 1457.800963
$ws.Cells.Item(11, 15).Value = 0.4649949779817998
$ws.Cells.Item(11, 16).Value = 0.4724755698027463
$ws.Cells.Item(11, 17).Value = 393.4332676290574
$ws.Cells.Item(11, 18).Value = 3540.899408661517
$ws.Cells.Item(11, 19).Value = 0.3629329536373718
$ws.Cells.Item(11, 20).Value = 0.3687716259093075
$ws.Cells.Item(12, 4).Value = "ECs"
$ws.Cells.Item(12, 5).Value = 1.0
$ws.Cells.Item(12, 6).Value = 0.3333333333333333
$ws.Cells.Item(12, 7).Value = 0.036349
$ws.Cells.Item(12, 8).Value = 0.109047
$ws.Cells.Item(12, 9).Value = 0.0350410011879885
$ws.Cells.Item(12, 10).Value = 0.0350410011879885
$ws.Cells.Item(12, 11).Value = 3.0
$ws.Cells.Item(12, 13).Value = 47.50824466666666
$ws.Cells.Item(12, 14).Value = 142.524734
$ws.Cells.Item(12, 15).Value = 0.04546113442798697
$ws.Cells.Item(12, 16).Value = 0.04619248897260801
$ws.Cells.Item(12, 17).Value = 1.726877185388666
$ws.Cells.Item(12, 18).Value = 15.541894668498
$ws.Cells.Item(12, 19).Value = 0.001593003665498396
$ws.Cells.Item(12, 20).Value = 0.001618631060965303
$ws.Cells.Item(13, 4).Value = "FAPs"
$ws.Cells.Item(13, 5).Value = 1.0
$ws.Cells.Item(13, 6).Value = 0.3333333333333333
$ws.Cells.Item(13, 7).Value = 0.036349
$ws.Cells.Item(13, 8).Value = 0.109047
$ws.Cells.Item(13, 9).Value = 0.0350410011879885
$ws.Cells.Item(13, 10).Value = 0.0350410011879885
$ws.Cells.Item(13, 11).Value = 3.0
$ws.Cells.Item(13, 15).Value = 0.105812544913079
$ws.Cells.Item(13, 16).Value = 0.1075148008416609
$ws.Cells.Item(13, 17).Value = 4.019373296276999
$ws.Cells.Item(13, 18).Value = 36.174359666493
$ws.Cells.Item(13, 19).Value = 0.003707777512003288
$ws.Cells.Item(13, 20).Value = 0.003767426264018986
$ws.Cells.Item(14, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(14, 5).Value = 1.0
$ws.Cells.Item(14, 6).Value = 0.3333333333333333
$ws.Cells.Item(14, 7).Value = 0.036349
$ws.Cells.Item(14, 8).Value = 0.109047
$ws.Cells.Item(14, 9).Value = 0.0350410011879885
$ws.Cells.Item(14, 10).Value = 0.0350410011879885
$ws.Cells.Item(14, 11).Value = 3.0
$ws.Cells.Item(14, 13).Value = 351.3736063333333
$ws.Cells.Item(14, 14).Value = 1054.120819
$ws.Cells.Item(14, 15).Value = 0.3362330657350935
$ws.Cells.Item(14, 16).Value = 0.3416422044152282
$ws.Cells.Item(14, 17).Value = 12.77207921661033
$ws.Cells.Item(14, 18).Value = 114.948712949493
$ws.Cells.Item(14, 19).Value = 0.01178194325586443
$ws.Cells.Item(14, 20).Value = 0.01197148489078102
$ws.Cells.Item(15, 4).Value = "MuSCs"
$ws.Cells.Item(15, 5).Value = 1.0
$ws.Cells.Item(15, 6).Value = 0.3333333333333333
$ws.Cells.Item(15, 7).Value = 0.036349
$ws.Cells.Item(15, 8).Value = 0.109047
$ws.Cells.Item(15, 9).Value = 0.0350410011879885
$ws.Cells.Item(15, 10).Value = 0.0350410011879885
$ws.Cells.Item(15, 11).Value = 2.0
$ws.Cells.Item(15, 13).Value = 49.6371195
$ws.Cells.Item(15, 14).Value = 99.274239
$ws.Cells.Item(15, 15).Value = 0.04749827694204053
$ws.Cells.Item(15, 16).Value = 0.03217493596775667
$ws.Cells.Item(15, 17).Value = 1.8042596567055
$ws.Cells.Item(15, 18).Value = 10.825557940233
$ws.Cells.Item(15, 19).Value = 0.001664387178753449
$ws.Cells.Item(15, 20).Value = 0.001127441969469615
$ws.Cells.Item(16, 5).Value = 1.0
$ws.Cells.Item(16, 6).Value = 0.3333333333333333
$ws.Cells.Item(16, 7).Value = 0.036349
$ws.Cells.Item(16, 8).Value = 0.109047
$ws.Cells.Item(16, 9).Value = 0.0350410011879885
$ws.Cells.Item(16, 10).Value = 0.0350410011879885
$ws.Cells.Item(16, 11).Value = 3.0
$ws.Cells.Item(16, 13).Value = 485.9336543333334
$ws.Cells.Item(16, 14).Value = 1457.800963
$ws.Cells.Item(16, 15).Value = 0.4649949779817998
$ws.Cells.Item(16, 16).Value = 0.4724755698027463
$ws.Cells.Item(16, 17).Value = 17.66320240136233
$ws.Cells.Item(16, 18).Value = 158.968821612261
$ws.Cells.Item(16, 19).Value = 0.01629388957586893
$ws.Cells.Item(16, 20).Value = 0.01655601700275357
$ws.Cells.Item(17, 4).Value = "ECs"
$ws.Cells.Item(17, 5).Value = 1.0
$ws.Cells.Item(17, 7).Value = 0.141536
$ws.Cells.Item(17, 8).Value = 0.424608
$ws.Cells.Item(17, 9).Value = 0.1364429047330914
$ws.Cells.Item(17, 10).Value = 0.1364429047330914
$ws.Cells.Item(17, 11).Value = 3.0
$ws.Cells.Item(17, 13).Value = 47.50824466666666
$ws.Cells.Item(17, 14).Value = 142.524734
$ws.Cells.Item(17, 15).Value = 0.04546113442798697
$ws.Cells.Item(17, 16).Value = 0.04619248897260801
$ws.Cells.Item(17, 17).Value = 6.724126917141333
$ws.Cells.Item(17, 18).Value = 60.517142254272
$ws.Cells.Item(17, 19).Value = 0.006202849233816087
$ws.Cells.Item(17, 20).Value = 0.00630263737227393
$ws.Cells.Item(18, 4).Value = "FAPs"
$ws.Cells.Item(18, 5).Value = 1.0
$ws.Cells.Item(18, 7).Value = 0.141536
$ws.Cells.Item(18, 8).Value = 0.424608
$ws.Cells.Item(18, 9).Value = 0.1364429047330914
$ws.Cells.Item(18, 10).Value = 0.1364429047330914
$ws.Cells.Item(18, 11).Value = 3.0
$ws.Cells.Item(18, 15).Value = 0.105812544913079
$ws.Cells.Item(18, 16).Value = 0.1075148008416609
$ws.Cells.Item(18, 17).Value = 15.650664911328
$ws.Cells.Item(18, 18).Value = 140.855984201952
$ws.Cells.Item(18, 19).Value = 0.0144373709851412
$ws.Cells.Item(18, 20).Value = 0.01466963172863603
$ws.Cells.Item(19, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(19, 5).Value = 1.0
$ws.Cells.Item(19, 7).Value = 0.141536
$ws.Cells.Item(19, 8).Value = 0.424608
$ws.Cells.Item(19, 9).Value = 0.1364429047330914
$ws.Cells.Item(19, 10).Value = 0.1364429047330914
$ws.Cells.Item(19, 11).Value = 3.0
$ws.Cells.Item(19, 13).Value = 351.3736063333333
$ws.Cells.Item(19, 14).Value = 1054.120819
$ws.Cells.Item(19, 15).Value = 0.3362330657350935
$ws.Cells.Item(19, 16).Value = 0.3416422044152282
$ws.Cells.Item(19, 17).Value = 49.73201474599466
$ws.Cells.Item(19, 18).Value = 447.5881327139519
$ws.Cells.Item(19, 19).Value = 0.04587661615620862
$ws.Cells.Item(19, 20).Value = 0.04661465474983031
$ws.Cells.Item(20, 4).Value = "MuSCs"
$ws.Cells.Item(20, 5).Value = 1.0
$ws.Cells.Item(20, 7).Value = 0.141536
$ws.Cells.Item(20, 8).Value = 0.424608
$ws.Cells.Item(20, 9).Value = 0.1364429047330914
$ws.Cells.Item(20, 10).Value = 0.1364429047330914
$ws.Cells.Item(20, 11).Value = 2.0
$ws.Cells.Item(20, 13).Value = 49.6371195
$ws.Cells.Item(20, 14).Value = 99.274239
$ws.Cells.Item(20, 15).Value = 0.04749827694204053
$ws.Cells.Item(20, 16).Value = 0.03217493596775667
$ws.Cells.Item(20, 17).Value = 7.025439345551999
$ws.Cells.Item(20, 18).Value = 42.152636073312
$ws.Cells.Item(20, 19).Value = 0.006480802875788828
$ws.Cells.Item(20, 20).Value = 0.004390041723041939
$ws.Cells.Item(21, 7).Value = 0.141536
$ws.Cells.Item(21, 8).Value = 0.424608
$ws.Cells.Item(21, 9).Value = 0.1364429047330914
$ws.Cells.Item(21, 10).Value = 0.1364429047330914
$ws.Cells.Item(21, 11).Value = 3.0
$ws.Cells.Item(21, 13).Value = 485.9336543333334
$ws.Cells.Item(21, 14).Value = 1457.800963
$ws.Cells.Item(21, 15).Value = 0.4649949779817998
$ws.Cells.Item(21, 16).Value = 0.4724755698027463
$ws.Cells.Item(21, 17).Value = 68.77710569972267
$ws.Cells.Item(21, 18).Value = 618.993951297504
$ws.Cells.Item(21, 19).Value = 0.06344526548213664
$ws.Cells.Item(21, 20).Value = 0.06446593915930919
